# Apply the authored edit: add a new employee row (SHAN / MANAGER) to the
# EMP sheet, format the SAL/COMM columns as Currency, and make EMP the
# active sheet/tab again (it had drifted to ExtraSheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EMP")

# --- format SAL (F) and COMM (G) as Currency, rows 2:14 -------------------
$ws.Range("F2:G13").Style = "Currency"

# --- new row 14: EMPNO 7936, SHAN, MANAGER, MGR 7839, same hiredate as
#     row 13, SAL with a floating-point value, no COMM, DEPTNO 10 ----------
$ws.Range("A14").Value = 7936
$ws.Range("B14").Value = "SHAN"
$ws.Range("C14").Value = "MANAGER"
$ws.Range("D14").Value = 7839
$ws.Range("E14").Value = $ws.Range("E13").Value2
$ws.Range("F14").Value = 1300.9654
$ws.Range("F14").Style = "Currency"
$ws.Range("G14").Style = "Currency"
$ws.Range("H14").Value = 10

# --- column widths: let Excel re-fit F/G now that they hold currency text -
$ws.Columns("F:G").AutoFit()

# --- restore EMP as the active sheet / selected cell -----------------------
$ws.Activate()
$ws.Range("F14").Select()
